$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, Year, Month, Day, Calidad (I), Volumen (J), PrecioMinimo (K), PrecioMaximo (L), PrecioPromedio (M), PrecioKg (P)
$rows = @(
    @(2, 2021, 3, 25, "Primera", 200, 10000, 11000, 10500, 583),
    @(3, 2021, 3, 25, "Segunda", 100, 9000, 9000, 9000, 500),
    @(4, 2021, 3, 17, "Primera", 100, 10000, 11000, 10500, 583),
    @(5, 2021, 3, 17, "Segunda", 50, 9000, 9000, 9000, 500),
    @(6, 2021, 6, 3, "Primera", 100, 13000, 14000, 13500, 750),
    @(7, 2021, 6, 3, "Segunda", 50, 11000, 11000, 11000, 611),
    @(8, 2021, 5, 4, "Primera", 100, 9000, 10000, 9500, 528),
    @(9, 2021, 5, 4, "Segunda", 50, 8000, 8000, 8000, 444),
    @(10, 2021, 4, 22, "Primera", 200, 10000, 11000, 10500, 583),
    @(11, 2021, 4, 22, "Segunda", 50, 8000, 8000, 8000, 444),
    @(12, 2021, 3, 4, "Primera", 100, 12000, 13000, 12500, 694),
    @(13, 2021, 3, 4, "Segunda", 50, 10000, 10000, 10000, 556),
    @(14, 2021, 7, 20, "Primera", 100, 14000, 15000, 14500, 806),
    @(15, 2021, 2, 26, "Primera", 100, 12000, 12000, 12000, 667),
    @(16, 2021, 2, 26, "Segunda", 100, 10000, 10000, 10000, 556),
    @(17, 2021, 4, 30, "Primera", 200, 10000, 11000, 10500, 583),
    @(18, 2021, 4, 30, "Segunda", 100, 9000, 9000, 9000, 500),
    @(19, 2021, 6, 9, "Primera", 100, 10000, 11000, 10500, 583),
    @(20, 2021, 6, 9, "Segunda", 50, 9000, 9000, 9000, 500),
    @(21, 2021, 5, 26, "Primera", 100, 11000, 12000, 11500, 639),
    @(22, 2021, 5, 26, "Segunda", 50, 9000, 9000, 9000, 500),
    @(23, 2021, 7, 2, "Primera", 200, 15000, 16000, 15500, 861),
    @(24, 2021, 7, 2, "Segunda", 100, 13000, 13000, 13000, 722),
    @(25, 2021, 5, 12, "Primera", 200, 9000, 10000, 9500, 528),
    @(26, 2021, 5, 12, "Segunda", 100, 8000, 8000, 8000, 444),
    @(27, 2021, 5, 19, "Primera", 100, 12000, 13000, 12500, 694),
    @(28, 2021, 5, 19, "Segunda", 50, 10000, 10000, 10000, 556),
    @(29, 2021, 7, 15, "Primera", 200, 15000, 16000, 15500, 861),
    @(30, 2021, 7, 15, "Segunda", 100, 14000, 14000, 14000, 778),
    @(31, 2021, 3, 10, "Primera", 100, 13000, 14000, 13500, 750),
    @(32, 2021, 4, 27, "Primera", 200, 10000, 11000, 10500, 583),
    @(33, 2021, 4, 27, "Segunda", 100, 9000, 9000, 9000, 500),
    @(34, 2022, 2, 23, "Primera", 100, 15000, 16000, 15500, 861),
    @(35, 2021, 6, 2, "Primera", 200, 11000, 12000, 11500, 639),
    @(36, 2021, 6, 2, "Segunda", 100, 10000, 10000, 10000, 556),
    @(37, 2021, 6, 17, "Primera", 100, 13000, 14000, 13500, 750),
    @(38, 2021, 6, 17, "Segunda", 50, 11000, 11000, 11000, 611),
    @(39, 2021, 7, 7, "Primera", 200, 15000, 16000, 15500, 861),
    @(40, 2021, 7, 7, "Segunda", 50, 13000, 13000, 13000, 722)
)

foreach ($row in $rows) {
    $r = $row[0]
    $y = $row[1]
    $mo = $row[2]
    $da = $row[3]
    $calidad = $row[4]
    $volumen = $row[5]
    $precioMin = $row[6]
    $precioMax = $row[7]
    $precioProm = $row[8]
    $precioKg = $row[9]

    $ws.Cells.Item($r, 4).Value = (Get-Date -Year $y -Month $mo -Day $da -Hour 0 -Minute 0 -Second 0)
    $ws.Cells.Item($r, 9).Value = $calidad
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $precioMin
    $ws.Cells.Item($r, 12).Value = $precioMax
    $ws.Cells.Item($r, 13).Value = $precioProm
    $ws.Cells.Item($r, 16).Value = $precioKg
}

"Done updating rows 2-40."
